$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in the header (F1).
$ws.Range("F1").Value = "Last status check on: 04.02.2022 13:30"

# Tesco row (row 3): the price delta and the "Old Datum" timestamp were
# previously stored as text; refresh them as real numeric/date values,
# matching how the other rows in the sheet are stored.
$ws.Range("D3").Value = -0.01

$ws.Range("E3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 44596.55221064815
